# 125I Radioactive Disposal Log & Printout
# The sheet is converted from a 3H-Ligand sink/dry-waste log (in uCi) into a
# 125I-Ligand batch log (in mCi), and a "Name" column is added at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-label the existing header cells for the new 125I workflow ---
$ws.Range("B1").Value = "125I-Ligand"
$ws.Range("C1").Value = "Batch Number"
$ws.Range("D1").Value = "mCi"
$ws.Range("E1").Value = "Liquid Disposal (mCi)"
$ws.Range("F1").Value = "Dry Waste (mCi)"

# --- "Name" header moves out to the new column G, keeping the bold/centered
#     header style used by the rest of row 1 ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$ws.Range("G1").Value = "Name"

# --- Column widths, re-fit for the new header text ---
$ws.Columns.Item(2).ColumnWidth = 10.084
$ws.Columns.Item(3).ColumnWidth = 12.7507
$ws.Columns.Item(4).ColumnWidth = 11.584
$ws.Columns.Item(5).ColumnWidth = 18.9173
$ws.Columns.Item(6).ColumnWidth = 14.084

# --- Selection moves to F2 ---
$ws.Range("F2").Select() | Out-Null
